$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("E5").Value = 1.485473821631866
$ws.Range("C6").Value = 1.21254482274098
$ws.Range("E6").Value = 1.799394172339364
$ws.Range("C8").Value = 0.4712609263772816
$ws.Range("E8").Value = 1.107727073902165
$ws.Range("E10").Value = 2.221748592150141
$ws.Range("C11").Value = 4.109890522944326
$ws.Range("E13").Value = 1.649865498505254
$ws.Range("E14").Value = 2.284828905445191
$ws.Range("C16").Value = 2.777797690741446
$ws.Range("E17").Value = 1.804067895915296
$ws.Range("C18").Value = -1.432689847121826
$ws.Range("C19").Value = 2.033479419175155
$ws.Range("E19").Value = 2.146365108912263
